# Commit: "added connection to SSRS"
# The "Analysis" worksheet's "Reported Volume (uL)" column (C) held
# placeholder/hardcoded 100.0 values for every well. This updates those
# cells with the real reported volumes now pulled in from the SSRS
# connection (one value per well row; the handful of rows not listed in
# the source diff keep their existing 100.0 placeholder).
# Dependent formula cells (D "Result", E "Difference (uL)", F
# "Difference (%)") recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Analysis")

$ws.Range("C4").Value = 56.0
$ws.Range("C8").Value = 156.0441
$ws.Range("C9").Value = 124.3936
$ws.Range("C10").Value = 123.3568
$ws.Range("C11").Value = 135.9624
$ws.Range("C14").Value = 127.5859
$ws.Range("C15").Value = 104.6938
$ws.Range("C16").Value = 137.9815
$ws.Range("C17").Value = 145.9214
$ws.Range("C18").Value = 152.961
$ws.Range("C19").Value = 141.2284
$ws.Range("C20").Value = 156.9991
$ws.Range("C21").Value = 141.8833
$ws.Range("C22").Value = 138.2271
$ws.Range("C23").Value = 166.9309
$ws.Range("C24").Value = 143.4385
$ws.Range("C25").Value = 127.4767
$ws.Range("C26").Value = 135.4986
$ws.Range("C27").Value = 134.6254
$ws.Range("C28").Value = 128.3226
$ws.Range("C29").Value = 151.0783
$ws.Range("C30").Value = 139.6187
$ws.Range("C31").Value = 148.8682
$ws.Range("C32").Value = 131.515
$ws.Range("C33").Value = 140.9556
$ws.Range("C34").Value = 107.4496
$ws.Range("C35").Value = 147.2311
$ws.Range("C36").Value = 150.9965
$ws.Range("C37").Value = 123.2749
$ws.Range("C38").Value = 185.1845
$ws.Range("C39").Value = 153.3702
$ws.Range("C40").Value = 144.4208
$ws.Range("C41").Value = 118.1453
$ws.Range("C42").Value = 135.1166
$ws.Range("C43").Value = 143.4658
$ws.Range("C44").Value = 134.5163
$ws.Range("C45").Value = 134.2434
$ws.Range("C46").Value = 170.6963
$ws.Range("C47").Value = 134.3799
$ws.Range("C48").Value = 154.1342
$ws.Range("C49").Value = 155.5803
$ws.Range("C50").Value = 144.5026
$ws.Range("C51").Value = 142.5381
$ws.Range("C52").Value = 151.9242
$ws.Range("C53").Value = 149.6322
$ws.Range("C54").Value = 133.9434
$ws.Range("C55").Value = 135.9351
$ws.Range("C56").Value = 139.0183
$ws.Range("C57").Value = 148.3226
$ws.Range("C58").Value = 140.6282
$ws.Range("C59").Value = 153.0974
$ws.Range("C60").Value = 143.7932
$ws.Range("C61").Value = 154.2707
$ws.Range("C62").Value = 134.9801
$ws.Range("C63").Value = 173.2065
$ws.Range("C64").Value = 148.3226
$ws.Range("C65").Value = 144.6936
$ws.Range("C66").Value = 176.5352
$ws.Range("C67").Value = 158.0633
$ws.Range("C68").Value = 138.3635
$ws.Range("C69").Value = 154.6799
$ws.Range("C70").Value = 141.6923
$ws.Range("C71").Value = 140.1916
$ws.Range("C72").Value = 126.6037
$ws.Range("C73").Value = 140.3826
$ws.Range("C74").Value = 140.0006
$ws.Range("C75").Value = 146.3307
$ws.Range("C76").Value = 128.7046
$ws.Range("C77").Value = 131.133
$ws.Range("C78").Value = 138.5
$ws.Range("C79").Value = 139.5913
$ws.Range("C80").Value = 115.6624
$ws.Range("C81").Value = 156.699
$ws.Range("C82").Value = 144.639
$ws.Range("C83").Value = 157.5448
$ws.Range("C84").Value = 0.0
$ws.Range("C85").Value = 0.0
$ws.Range("C86").Value = 0.0
$ws.Range("C87").Value = 0.0
$ws.Range("C88").Value = 0.0
$ws.Range("C89").Value = 0.0
$ws.Range("C90").Value = 0.0
$ws.Range("C91").Value = 0.0
$ws.Range("C92").Value = 0.0
$ws.Range("C93").Value = 0.0
$ws.Range("C94").Value = 0.0
$ws.Range("C95").Value = 0.0
$ws.Range("C96").Value = 0.0
$ws.Range("C97").Value = 0.0
$ws.Range("C98").Value = 0.0
$ws.Range("C99").Value = 0.0
